$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2637.6667
$ws.Range("I129").Value = 6790.6875
$ws.Range("J129").Value = 1016.9756
$ws.Range("K129").Value = 20372.0625
$ws.Range("L129").Value = 3050.9268
$ws.Range("M129").Value = -15372.0625
$ws.Range("N129").Value = -13050.9268

$ws.Range("H137").Value = 1576.0385
$ws.Range("I137").Value = 1624
$ws.Range("J137").Value = 1499.3
$ws.Range("K137").Value = 4872
$ws.Range("L137").Value = 4497.9
$ws.Range("M137").Value = -2322
$ws.Range("N137").Value = -9597.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H45").Value = 127153.25
$ws.Range("I45").Value = 167877.67
$ws.Range("K45").Value = 167877.67
$ws.Range("M45").Value = -167500.67

$ws.Range("H61").Value = 1363.5428
$ws.Range("I61").Value = 1134.4445
$ws.Range("J61").Value = 2136.75
$ws.Range("K61").Value = 1134.4445
$ws.Range("L61").Value = 2136.75
$ws.Range("M61").Value = -922.4445000000001
$ws.Range("N61").Value = -2560.75

$ws.Range("H74").Value = 698.13043
$ws.Range("I74").Value = 654
$ws.Range("J74").Value = 823.1667
$ws.Range("K74").Value = 654
$ws.Range("L74").Value = 823.1667
$ws.Range("M74").Value = 220
$ws.Range("N74").Value = -2571.1667

$ws.Range("H77").Value = 698.13043
$ws.Range("I77").Value = 654
$ws.Range("J77").Value = 823.1667
$ws.Range("K77").Value = 3270
$ws.Range("L77").Value = 4115.8335
$ws.Range("M77").Value = 1098
$ws.Range("N77").Value = -12851.8335

$ws.Range("H122").Value = 958
$ws.Range("I122").Value = 987.25
$ws.Range("J122").Value = 899.5
$ws.Range("K122").Value = 2961.75
$ws.Range("L122").Value = 2698.5
$ws.Range("M122").Value = -511.75
$ws.Range("N122").Value = -7598.5

$ws.Range("H132").Value = 1382.6885
$ws.Range("I132").Value = 1106.9636
$ws.Range("J132").Value = 3910.1667
$ws.Range("K132").Value = 3320.8908
$ws.Range("L132").Value = 11730.5001
$ws.Range("M132").Value = -790.8908000000001
$ws.Range("N132").Value = -16790.5001

$ws.Range("H136").Value = 1363.5428
$ws.Range("I136").Value = 1134.4445
$ws.Range("J136").Value = 2136.75
$ws.Range("K136").Value = 3403.3335
$ws.Range("L136").Value = 6410.25
$ws.Range("M136").Value = -853.3335000000002
$ws.Range("N136").Value = -11510.25

$ws.Range("H137").Value = 43000
$ws.Range("J137").Value = 43000
$ws.Range("L137").Value = 43000
$ws.Range("N137").Value = -53200

$ws.Range("H138").Value = 53475.6
$ws.Range("J138").Value = 53475.6
$ws.Range("L138").Value = 53475.6
$ws.Range("N138").Value = -63755.6

$ws.Range("H139").Value = 32619.166
$ws.Range("J139").Value = 35143
$ws.Range("L139").Value = 35143
$ws.Range("N139").Value = -45423

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3291.7144
$ws.Range("I134").Value = 2904.1082
$ws.Range("J134").Value = 6160
$ws.Range("K134").Value = 8712.3246
$ws.Range("L134").Value = 18480
$ws.Range("M134").Value = -6177.3246
$ws.Range("N134").Value = -23550

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2349.75
$ws.Range("I62").Value = 2099.7144
$ws.Range("K62").Value = 2099.7144
$ws.Range("M62").Value = -1475.7144

$ws.Range("H65").Value = 2349.75
$ws.Range("I65").Value = 2099.7144
$ws.Range("K65").Value = 10498.572
$ws.Range("M65").Value = -7378.572

$ws.Range("H132").Value = 4691.8335
$ws.Range("I132").Value = 5509.7
$ws.Range("J132").Value = 4107.643
$ws.Range("K132").Value = 16529.1
$ws.Range("L132").Value = 12322.929
$ws.Range("M132").Value = -13999.1
$ws.Range("N132").Value = -17382.929

$ws.Range("H134").Value = 2165.2222
$ws.Range("I134").Value = 1722
$ws.Range("J134").Value = 2519.8
$ws.Range("K134").Value = 5166
$ws.Range("L134").Value = 7559.400000000001
$ws.Range("M134").Value = -2631
$ws.Range("N134").Value = -12629.4

$ws.Range("H138").Value = 58499.215
$ws.Range("J138").Value = 58499.215
$ws.Range("L138").Value = 58499.215
$ws.Range("N138").Value = -68779.215

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 664.8889
$ws.Range("I18").Value = 500.8
$ws.Range("K18").Value = 1502.4
$ws.Range("M18").Value = -1333.4

$ws.Range("H122").Value = 468.625

$ws.Range("H140").Value = 5190.0713
$ws.Range("I140").Value = 6865.8335
$ws.Range("J140").Value = 2173.7
$ws.Range("K140").Value = 20597.5005
$ws.Range("L140").Value = 6521.099999999999
$ws.Range("M140").Value = -15417.5005
$ws.Range("N140").Value = -16881.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 3803.4
$ws.Range("I33").Value = 3017
$ws.Range("K33").Value = 3017
$ws.Range("M33").Value = -2765

$ws.Range("H122").Value = 899.8
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 833
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 2499
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -7399

$ws.Range("H132").Value = 2491.2258
$ws.Range("I132").Value = 2181.4075
$ws.Range("J132").Value = 4582.5
$ws.Range("K132").Value = 6544.2225
$ws.Range("L132").Value = 13747.5
$ws.Range("M132").Value = -4014.2225
$ws.Range("N132").Value = -18807.5

$ws.Range("H133").Value = 42745
$ws.Range("J133").Value = 42745
$ws.Range("L133").Value = 42745
$ws.Range("N133").Value = -52865

$ws.Range("H139").Value = 58979
$ws.Range("J139").Value = 58979
$ws.Range("L139").Value = 58979
$ws.Range("N139").Value = -69259

$ws.Range("H140").Value = 108928.57
$ws.Range("J140").Value = 108928.57
$ws.Range("L140").Value = 108928.57
$ws.Range("N140").Value = -119288.57

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5673.067
$ws.Range("I7").Value = 5350.8
$ws.Range("J7").Value = 6317.6
$ws.Range("K7").Value = 5350.8
$ws.Range("L7").Value = 6317.6
$ws.Range("M7").Value = -5238.8
$ws.Range("N7").Value = -6541.6

$ws.Range("H64").Value = 21666.334
$ws.Range("J64").Value = 21666.334
$ws.Range("L64").Value = 21666.334
$ws.Range("N64").Value = -22116.334

$ws.Range("H67").Value = 21666.334
$ws.Range("J67").Value = 21666.334
$ws.Range("L67").Value = 21666.334
$ws.Range("N67").Value = -23226.334

$ws.Range("H126").Value = 5673.067
$ws.Range("I126").Value = 5350.8
$ws.Range("J126").Value = 6317.6
$ws.Range("K126").Value = 16052.4
$ws.Range("L126").Value = 18952.8
$ws.Range("M126").Value = -13582.4
$ws.Range("N126").Value = -23892.8

$ws.Range("H132").Value = 3424.7334
$ws.Range("I132").Value = 3334.926
$ws.Range("J132").Value = 4233
$ws.Range("K132").Value = 10004.778
$ws.Range("L132").Value = 12699
$ws.Range("M132").Value = -7474.778
$ws.Range("N132").Value = -17759

$ws.Range("H134").Value = 61666.668
$ws.Range("J134").Value = 61666.668
$ws.Range("L134").Value = 61666.668
$ws.Range("N134").Value = -71806.66800000001

$ws.Range("H136").Value = 1970.3182
$ws.Range("I136").Value = 1509.8667
$ws.Range("K136").Value = 4529.6001
$ws.Range("M136").Value = -1979.6001

$ws.Range("H140").Value = 88205.8
$ws.Range("J140").Value = 88205.8
$ws.Range("L140").Value = 88205.8
$ws.Range("N140").Value = -98565.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 39999.5
$ws.Range("J63").Value = 39999.5
$ws.Range("L63").Value = 39999.5
$ws.Range("N63").Value = -41247.5

$ws.Range("H66").Value = 39999.5
$ws.Range("J66").Value = 39999.5
$ws.Range("L66").Value = 119998.5
$ws.Range("N66").Value = -126238.5

$ws.Range("H122").Value = 1960.1111
$ws.Range("I122").Value = 1508.375
$ws.Range("J122").Value = 2321.5
$ws.Range("K122").Value = 4525.125
$ws.Range("L122").Value = 6964.5
$ws.Range("M122").Value = -2075.125
$ws.Range("N122").Value = -11864.5

$ws.Range("H126").Value = 1808
$ws.Range("I126").Value = 1623.2858
$ws.Range("J126").Value = 2066.6
$ws.Range("K126").Value = 4869.857400000001
$ws.Range("L126").Value = 6199.799999999999
$ws.Range("M126").Value = -2399.857400000001
$ws.Range("N126").Value = -11139.8

$ws.Range("H132").Value = 2672.2222
$ws.Range("I132").Value = 2535.353
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 7606.059
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -5076.059
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 695.5517
$ws.Range("I136").Value = 526.1111
$ws.Range("K136").Value = 1578.3333
$ws.Range("M136").Value = 971.6667000000002

$ws.Range("H138").Value = 66172.5
$ws.Range("J138").Value = 66172.5
$ws.Range("L138").Value = 66172.5
$ws.Range("N138").Value = -76452.5
